$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 33266.668
$ws.Range("J75").Value = 33266.668
$ws.Range("L75").Value = 33266.668
$ws.Range("N75").Value = -35138.668

$ws.Range("H78").Value = 33266.668
$ws.Range("J78").Value = 33266.668
$ws.Range("L78").Value = 99800.00399999999
$ws.Range("N78").Value = -109160.004

$ws.Range("H115").Value = 1135
$ws.Range("I115").Value = 1203.5
$ws.Range("J115").Value = 450
$ws.Range("K115").Value = 3610.5
$ws.Range("L115").Value = 1350
$ws.Range("M115").Value = -2043.5
$ws.Range("N115").Value = -4484

$ws.Range("H118").Value = 775.2
$ws.Range("I118").Value = 490
$ws.Range("J118").Value = 846.5
$ws.Range("K118").Value = 1470
$ws.Range("L118").Value = 2539.5
$ws.Range("M118").Value = 187
$ws.Range("N118").Value = -5853.5

$ws.Range("H137").Value = 1490139.8
$ws.Range("I137").Value = 2166194
$ws.Range("J137").Value = 2820.2
$ws.Range("K137").Value = 6498582
$ws.Range("L137").Value = 8460.599999999999
$ws.Range("M137").Value = -6496032
$ws.Range("N137").Value = -13560.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 9526852
$ws.Range("I88").Value = 16668741
$ws.Range("J88").Value = 4333.3335
$ws.Range("K88").Value = 16668741
$ws.Range("L88").Value = 4333.3335
$ws.Range("M88").Value = -16668335
$ws.Range("N88").Value = -5145.3335

$ws.Range("H91").Value = 9526852
$ws.Range("I91").Value = 16668741
$ws.Range("J91").Value = 4333.3335
$ws.Range("K91").Value = 16668741
$ws.Range("L91").Value = 4333.3335
$ws.Range("M91").Value = -16667337
$ws.Range("N91").Value = -7141.3335

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1949.5555
$ws.Range("I86").Value = 1566
$ws.Range("J86").Value = 2716.6667
$ws.Range("K86").Value = 1566
$ws.Range("L86").Value = 2716.6667
$ws.Range("M86").Value = -443
$ws.Range("N86").Value = -4962.6667

$ws.Range("H89").Value = 1949.5555
$ws.Range("I89").Value = 1566
$ws.Range("J89").Value = 2716.6667
$ws.Range("K89").Value = 7830
$ws.Range("L89").Value = 13583.3335
$ws.Range("M89").Value = -2214
$ws.Range("N89").Value = -24815.3335

$ws.Range("H134").Value = 2195.762
$ws.Range("I134").Value = 1555.55
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 4666.65
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -2131.65
$ws.Range("N134").Value = -50070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2800.4138
$ws.Range("I31").Value = 1054.4546
$ws.Range("K31").Value = 1054.4546
$ws.Range("M31").Value = -759.4546

$ws.Range("H34").Value = 2800.4138
$ws.Range("I34").Value = 1054.4546
$ws.Range("K34").Value = 1054.4546
$ws.Range("M34").Value = -852.4546

$ws.Range("H106").Value = 37500
$ws.Range("J106").Value = 37500
$ws.Range("L106").Value = 37500
$ws.Range("N106").Value = -40024

$ws.Range("H135").Value = 39926.668
$ws.Range("J135").Value = 39926.668
$ws.Range("L135").Value = 39926.668
$ws.Range("N135").Value = -50066.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 535151.25
$ws.Range("I5").Value = 599.4167
$ws.Range("J5").Value = 1028583.7
$ws.Range("K5").Value = 1798.2501
$ws.Range("L5").Value = 3085751.1
$ws.Range("M5").Value = -1686.2501
$ws.Range("N5").Value = -3085975.1

$ws.Range("H107").Value = 40348.44
$ws.Range("I107").Value = 338.9
$ws.Range("J107").Value = 67021.47
$ws.Range("K107").Value = 1016.7
$ws.Range("L107").Value = 201064.41
$ws.Range("M107").Value = 903.3000000000001
$ws.Range("N107").Value = -204904.41

$ws.Range("H131").Value = 660.8163500000001
$ws.Range("I131").Value = 243.12
$ws.Range("J131").Value = 803.86304
$ws.Range("K131").Value = 729.36
$ws.Range("L131").Value = 2411.58912
$ws.Range("M131").Value = 4310.64
$ws.Range("N131").Value = -12491.58912

$ws.Range("H132").Value = 2423.2258
$ws.Range("I132").Value = 970
$ws.Range("J132").Value = 3115.238
$ws.Range("K132").Value = 8730
$ws.Range("L132").Value = 28037.142
$ws.Range("M132").Value = -6200
$ws.Range("N132").Value = -33097.142

$ws.Range("H133").Value = 5000
$ws.Range("I133").Value = 5000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 15000
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -9940
$ws.Range("N133").Value = $null

$ws.Range("H134").Value = 3797.9644
$ws.Range("I134").Value = 3185.7222
$ws.Range("J134").Value = 4900
$ws.Range("K134").Value = 9557.1666
$ws.Range("L134").Value = 14700
$ws.Range("M134").Value = -4487.1666
$ws.Range("N134").Value = -24840

$ws.Range("H135").Value = 535151.25
$ws.Range("I135").Value = 599.4167
$ws.Range("J135").Value = 1028583.7
$ws.Range("K135").Value = 5394.7503
$ws.Range("L135").Value = 9257253.299999999
$ws.Range("M135").Value = -2859.7503
$ws.Range("N135").Value = -9262323.299999999

$ws.Range("H136").Value = 6657.5
$ws.Range("I136").Value = 6657.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 19972.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -14872.5
$ws.Range("N136").Value = $null

$ws.Range("H137").Value = 1757.5
$ws.Range("I137").Value = 843.3333
$ws.Range("J137").Value = 4500
$ws.Range("K137").Value = 2529.9999
$ws.Range("L137").Value = 13500
$ws.Range("M137").Value = 2570.0001
$ws.Range("N137").Value = -23700

$ws.Range("H139").Value = 1115.625
$ws.Range("I139").Value = 990
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 2970
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 2170
$ws.Range("N139").Value = -19280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 30780.2
$ws.Range("J80").Value = 30780.2
$ws.Range("L80").Value = 30780.2
$ws.Range("N80").Value = -32776.2

$ws.Range("H83").Value = 30780.2
$ws.Range("J83").Value = 30780.2
$ws.Range("L83").Value = 92340.60000000001
$ws.Range("N83").Value = -102324.6

$ws.Range("H132").Value = 4913.4165
$ws.Range("I132").Value = 4685.6895
$ws.Range("K132").Value = 14057.0685
$ws.Range("M132").Value = -11527.0685
